$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old "Tipo"/"single" to column E)
$ws.Columns.Item(4).Insert()

# New header for column D
$ws.Range("D1").Value = "MAE"

# Updated numeric values
$ws.Range("B2").Value = 0.5379370305603909
$ws.Range("C2").Value = 0.989288531347028
$ws.Range("D2").Value = 0.6137271608150739
